# Mise à jour de l'application
# Adds a new player row (Theo Owono) to the roster sheet, extends the
# conditional formatting range on column A to cover the new row, and moves
# the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New player data for row 30
$ws.Range("A30").Value = "Theo Owono"
# -4108 = xlCenter
$ws.Range("A30").HorizontalAlignment = -4108
$ws.Range("A30").VerticalAlignment = -4108
$ws.Range("B30").Value = 23
$ws.Range("C30").Value = (Get-Date -Year 2002 -Month 9 -Day 30).Date
$ws.Range("D30").Value = "MC"
$ws.Range("E30").Value = "1m87"

# Extend the conditional formatting applied to column A so it also covers
# the newly added row 30 (was A22:A29, now A22:A30) without disturbing the
# existing rules/dxf references.
$fc = $ws.Range("A22:A29").FormatConditions
for ($i = 1; $i -le $fc.Count; $i++) {
    $fc.Item($i).ModifyAppliesToRange($ws.Range("A22:A30"))
}

# Move the active selection to F28, matching the updated view state
$ws.Range("F28").Select()
